$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the "cryptos" list: update prices / 1h volume % for (almost) every
# coin, and replace a few rows whose coin entirely changed position/identity
# (PEPE<->USDe swap rows 30/31, Hedera<->Stellar swap rows 43/44, and the
# former SuiNetwork row 51 which is now Maker).
# Every value is written through .Formula with a leading "'" so Excel always
# stores it as literal text (matches the original t="inlineStr" cells and
# keeps values like "55.244.61", "0.0970" or "0.0<sub>3</sub>0760" intact
# instead of being reinterpreted as numbers/dates).

$ws.Range("D2").Formula = "'55.244.61"
$ws.Range("E2").Formula = "'  -2.40%  "
$ws.Range("D3").Formula = "'2.366.01"
$ws.Range("E3").Formula = "'  -4.89%  "
$ws.Range("E4").Formula = "'  -0.03%  "
$ws.Range("D5").Formula = "'475.31"
$ws.Range("E5").Formula = "'  -3.22%  "
$ws.Range("D6").Formula = "'145.58"
$ws.Range("E6").Formula = "'  -1.04%  "
$ws.Range("D7").Formula = "'0.999"
$ws.Range("E7").Formula = "'  +0.44%  "
$ws.Range("D8").Formula = "'0.504"
$ws.Range("E8").Formula = "'  -2.25%  "
$ws.Range("D9").Formula = "'2.367.14"
$ws.Range("E9").Formula = "'  -5.42%  "
$ws.Range("D10").Formula = "'0.0970"
$ws.Range("E10").Formula = "'  -0.96%  "
$ws.Range("D11").Formula = "'5.40"
$ws.Range("E11").Formula = "'  -6.50%  "
$ws.Range("D12").Formula = "'0.321"
$ws.Range("E12").Formula = "'  -3.51%  "
$ws.Range("E13").Formula = "'  +0.77%  "
$ws.Range("D14").Formula = "'2.778.51"
$ws.Range("E14").Formula = "'  -4.68%  "
$ws.Range("D15").Formula = "'55.533.74"
$ws.Range("E15").Formula = "'  -1.96%  "
$ws.Range("D16").Formula = "'20.21"
$ws.Range("E16").Formula = "'  -5.61%  "
$ws.Range("D17").Formula = "'0.0000131"
$ws.Range("E17").Formula = "'  -4.73%  "
$ws.Range("D18").Formula = "'2.364.61"
$ws.Range("E18").Formula = "'  -5.23%  "
$ws.Range("D19").Formula = "'4.55"
$ws.Range("E19").Formula = "'  -0.62%  "
$ws.Range("D20").Formula = "'313.91"
$ws.Range("E20").Formula = "'  -2.46%  "
$ws.Range("D21").Formula = "'9.67"
$ws.Range("E21").Formula = "'  -5.18%  "
$ws.Range("D22").Formula = "'0.998"
$ws.Range("E22").Formula = "'  +0.11%  "
$ws.Range("E23").Formula = "'  -3.41%  "
$ws.Range("D24").Formula = "'56.62"
$ws.Range("E24").Formula = "'  -3.78%  "
$ws.Range("E25").Formula = "'  +0.38%  "
$ws.Range("E26").Formula = "'  -4.43%  "
$ws.Range("D27").Formula = "'0.154"
$ws.Range("E27").Formula = "'  -7.22%  "
$ws.Range("D28").Formula = "'2.471.97"
$ws.Range("E28").Formula = "'  -4.53%  "
$ws.Range("E29").Formula = "'  -6.93%  "
$ws.Range("B30").Formula = "'USDe"
$ws.Range("C30").Formula = "'https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D30").Formula = "'1.00"
$ws.Range("E30").Formula = "'  +0.12%  "
$ws.Range("B31").Formula = "'PEPE"
$ws.Range("C31").Formula = "'https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D31").Formula = "'0.0" + [char]0x2083 + "0760"
$ws.Range("E31").Formula = "'  -4.16%  "
$ws.Range("D32").Formula = "'148.27"
$ws.Range("E32").Formula = "'  -0.43%  "
$ws.Range("E33").Formula = "'  -1.78%  "
$ws.Range("D34").Formula = "'1.47"
$ws.Range("E34").Formula = "'  -2.56%  "
$ws.Range("D35").Formula = "'5.05"
$ws.Range("E35").Formula = "'  -3.38%  "
$ws.Range("D36").Formula = "'1.09"
$ws.Range("E36").Formula = "'  -5.25%  "
$ws.Range("D37").Formula = "'3.56"
$ws.Range("E37").Formula = "'  -5.07%  "
$ws.Range("D38").Formula = "'0.823"
$ws.Range("E38").Formula = "'  -4.99%  "
$ws.Range("D39").Formula = "'33.51"
$ws.Range("E39").Formula = "'  -2.28%  "
$ws.Range("E40").Formula = "'  +0.70%  "
$ws.Range("D41").Formula = "'1.33"
$ws.Range("E41").Formula = "'  -1.04%  "
$ws.Range("D42").Formula = "'3.37"
$ws.Range("E42").Formula = "'  -4.57%  "
$ws.Range("B43").Formula = "'Stellar"
$ws.Range("C43").Formula = "'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D43").Formula = "'0.0948"
$ws.Range("E43").Formula = "'  +3.41%  "
$ws.Range("B44").Formula = "'Hedera"
$ws.Range("C44").Formula = "'https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D44").Formula = "'0.0533"
$ws.Range("E44").Formula = "'  -4.84%  "
$ws.Range("D45").Formula = "'0.578"
$ws.Range("E45").Formula = "'  -5.68%  "
$ws.Range("D46").Formula = "'10.17"
$ws.Range("E46").Formula = "'  -0.16%  "
$ws.Range("D47").Formula = "'254.18"
$ws.Range("E47").Formula = "'  -2.67%  "
$ws.Range("D48").Formula = "'0.0222"
$ws.Range("E48").Formula = "'  -3.23%  "
$ws.Range("D49").Formula = "'4.47"
$ws.Range("E49").Formula = "'  -8.46%  "
$ws.Range("D50").Formula = "'16.93"
$ws.Range("E50").Formula = "'  -4.21%  "
$ws.Range("B51").Formula = "'Maker"
$ws.Range("C51").Formula = "'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D51").Formula = "'1.780.56"
$ws.Range("E51").Formula = "'  -6.27%  "
